# Grades.xlsx edit: finished grading ps 7 (Homework 7 -> column N),
# added midterm study guide 2 solutions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Homework 7 scores (column N) for graded students.
# Row 10 and row 11 are not part of this grading pass (no N value added).
$scores = @{
    3  = "=0"
    4  = "=0"
    5  = "=12/12"
    6  = "=12/12"
    7  = "=12/12"
    8  = "=10/12"
    9  = "=11/12"
    12 = "=12/12"
    13 = "=12/12"
    14 = "=10.5/12"
    15 = "=10.5/12"
    16 = "=12/12"
    17 = "=12/12"
    18 = "=12/12"
    19 = "=12/12"
}

foreach ($row in $scores.Keys) {
    $ws.Range("N$row").Formula = $scores[$row]
}

# Update the selected cell to reflect where editing left off.
$ws.Range("N10").Select()
